# Update the division-problem worksheet table: replace each exercise's
# text with the new problem while keeping run formatting (font/size) intact.
$d = $word.ActiveDocument
$t = $d.Tables(1)

# Row 1 (table row 1)
$t.Cell(1,1).Range.Text = "48÷4="
$t.Cell(1,2).Range.Text = "14÷3="
$t.Cell(1,3).Range.Text = "90÷9="
$t.Cell(1,4).Range.Text = "27÷2="
$t.Cell(1,5).Range.Text = "35÷2="

# Row 2 (table row 5 - rows 2-4 are blank spacer rows)
$t.Cell(5,1).Range.Text = "79÷2="
$t.Cell(5,2).Range.Text = "93÷8="
$t.Cell(5,3).Range.Text = "87÷2="
$t.Cell(5,4).Range.Text = "17÷4="
$t.Cell(5,5).Range.Text = "66÷2="

# Row 3 (table row 9)
$t.Cell(9,1).Range.Text = "69÷2="
$t.Cell(9,2).Range.Text = "65÷3="
$t.Cell(9,3).Range.Text = "83÷5="
$t.Cell(9,4).Range.Text = "78÷2="
$t.Cell(9,5).Range.Text = "14÷4="

# Row 4 (table row 13)
$t.Cell(13,1).Range.Text = "37÷9="
$t.Cell(13,2).Range.Text = "17÷5="
$t.Cell(13,3).Range.Text = "59÷8="
$t.Cell(13,4).Range.Text = "73÷4="
$t.Cell(13,5).Range.Text = "65÷6="

# Row 5 (table row 17)
$t.Cell(17,1).Range.Text = "33÷2="
$t.Cell(17,2).Range.Text = "22÷8="
$t.Cell(17,3).Range.Text = "82÷4="
$t.Cell(17,4).Range.Text = "98÷8="
$t.Cell(17,5).Range.Text = "24÷8="

Write-Output "done"
